$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update attendance tracking values from 0 to 1 for specific rows/columns
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("H4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 1

$ws.Range("H13").Value = 1

$ws.Range("H14").Value = 1

$ws.Range("H15").Value = 1

$ws.Range("H16").Value = 1

$ws.Range("H17").Value = 1

$ws.Range("H18").Value = 1
